# Emacs cheatsheet update: add "copy config in tmux" row, a new
# "search and replace" section, and a new "functions" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 11 (Navigation section) for the tmux copy
#    config shortcut. This shifts every following row down by one.
$ws.Rows(11).Insert()
$ws.Range("A11").Value = "Navigation"
$ws.Range("B11").Value = "C-x C-spc"
$ws.Range("C11").Value = "go to previous mark"

# 2. Append a new "search and replace" section after the existing data
#    (old row 54 is now row 55; row 56 stays blank, matching the gap
#    pattern already used elsewhere in the sheet).
$ws.Range("A57").Value = "search and replace"
$ws.Range("B57").Value = "C-x n n"
$ws.Range("C57").Value = "Narrow to region"

$ws.Range("A58").Value = "search and replace"
$ws.Range("B58").Value = "C-x n w"
$ws.Range("C58").Value = "widen"

# 3. Append a new bold "functions" section header (row 59 left blank)
#    followed by four single-column entries.
$ws.Range("A60").Value = "functions"
$ws.Range("A60").Font.Bold = $true
$ws.Range("A60").Font.Size = 14
$ws.Range("A60").Font.Name = "Aptos Narrow (Body)"
$ws.Rows(60).RowHeight = 19

$ws.Range("A61").Value = "goto line"
$ws.Range("A62").Value = "ag"
$ws.Range("A63").Value = "undo-tree-visualisation"
$ws.Range("A64").Value = "ace swap windows"

# 4. Restore the active selection to D1, matching the saved workbook
#    view state.
$ws.Range("D1").Select()
